$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.332.85"
$ws.Range("D3").Value = "3.679.86"
$ws.Range("E3").Value = "  -3.51%  "
$ws.Range("D5").Value = "682.79"
$ws.Range("D6").Value = "162.54"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("D7").Value = "3.678.51"
$ws.Range("E7").Value = "  -3.55%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -4.40%  "
$ws.Range("E10").Value = "  -7.53%  "
$ws.Range("D11").Value = "7.22"
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E13").Value = "  -6.16%  "
$ws.Range("D14").Value = "33.53"
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("D15").Value = "4.301.59"
$ws.Range("E15").Value = "  -3.50%  "
$ws.Range("D16").Value = "3.689.06"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "69.378.29"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "16.33"
$ws.Range("E19").Value = "  -5.52%  "
$ws.Range("E20").Value = "  -6.83%  "
$ws.Range("D21").Value = "478.04"
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("D22").Value = "9.82"
$ws.Range("E22").Value = "  -7.63%  "
$ws.Range("E23").Value = "  -8.20%  "
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("D25").Value = "3.826.19"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("E26").Value = "  -10.89%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  -7.49%  "
$ws.Range("D30").Value = "1.84"
$ws.Range("E30").Value = "  -10.29%  "
$ws.Range("E31").Value = "  -10.82%  "
$ws.Range("E32").Value = "  -5.88%  "
$ws.Range("D33").Value = "6.90"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").Value = "27.15"
$ws.Range("E34").Value = "  -6.82%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("D37").Value = "3.643.28"
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("E38").Value = "  -5.64%  "
$ws.Range("E40").Value = "  -7.24%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "0.958"
$ws.Range("E44").Value = "  -7.79%  "
$ws.Range("D45").Value = "48.23"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "154.33"
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("D47").Value = "2.85"
$ws.Range("E47").Value = "  -11.90%  "
$ws.Range("D48").Value = "1.33"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  -12.40%  "
$ws.Range("D50").Value = "393.14"
$ws.Range("E50").Value = "  -7.16%  "
$ws.Range("D51").Value = "8.10"
$ws.Range("E51").Value = "  -6.03%  "
